$d = $word.ActiveDocument

# --- Part 1: extend the first paragraph with the red "branch alternate" note ---
$p1 = $d.Paragraphs(1)
$rng = $p1.Range
$rng.MoveEnd(1, -1)       # exclude the paragraph mark so we stay inside paragraph 1
$rng.Collapse(0)          # wdCollapseEnd -> sit right after "document." before the mark
$rng.InsertAfter("  ")    # two trailing spaces appended to the existing run's text

$dash = [char]0x2013
$rng.Collapse(0)
$rng.InsertAfter("(This is a change " + $dash + " Ve")
$run1 = $d.Range($rng.Start, $rng.End)
$run1.Font.Color = 192            # 0x0000C0 -> w:val="C00000"

$rng.Collapse(0)
$rng.InsertAfter("rsion for branch alternate")
$run2 = $d.Range($rng.Start, $rng.End)
$run2.Font.Color = 192

$rng.Collapse(0)
$rng.InsertAfter(")")
$run3 = $d.Range($rng.Start, $rng.End)
$run3.Font.Color = 192

# --- Part 2: insert a new, empty, pre-formatted paragraph after "It will be treated..." ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(3)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newp.Range.InsertXML($xml)

Write-Output "edit complete"
